# "updated legacy GSC export data"
#
# The "Chart" sheet's first data row (2025-10-24) only ever had a partial
# export (no Not-indexed/Indexed counts), so it is dropped and every
# subsequent day shifts up by one row. The table now ends one day later
# (2026-01-20) and that day's Impressions count (previously missing) is
# filled in with the newly-available value of 44.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the obsolete 2025-10-24 row; Excel shifts rows 3:90 up to 2:89.
$ws.Rows.Item(2).Delete()

# The last row (now row 89, date 2026-01-20) gets its Impressions value.
$ws.Cells.Item(89, 4).Value = 44
